$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 - Machine 30
$ws.Range("A31").Value = 10030
$ws.Range("B31").Value = "Machine 30"
$ws.Range("C31").Value = "70-5A-0F-8C-01-39"
$ws.Range("D31").Value = "FB5962911664"
$ws.Range("E31").Value = "192.168.0.356"
$ws.Range("F31").Value = 1001
$ws.Range("G31").Value = "eng"
$ws.Range("H31").Value = $true
$ws.Range("I31").Value = "superadmin"
$ws.Range("J31").Value = "now()"

# Row 32 - Machine 31
$ws.Range("A32").Value = 10031
$ws.Range("B32").Value = "Machine 31"
$ws.Range("C32").Value = "58-20-B1-DA-F3-FB"
$ws.Range("D32").Value = "FB5962911663"
$ws.Range("E32").Value = "192.168.0.357"
$ws.Range("F32").Value = 1001
$ws.Range("G32").Value = "eng"
$ws.Range("H32").Value = $true
$ws.Range("I32").Value = "superadmin"
$ws.Range("J32").Value = "now()"

# Move the "end of data" selection marker down, same as Excel would after
# entering data in the last two rows (selects the rest of the column range).
[void]$ws.Range("A33:XFD1048576").Select()
